# Swap the contents of rows 2 and 3 on the active worksheet: the record
# that was in row 2 moves to row 3, and the record that was in row 3
# moves to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns used anywhere on the sheet (header row spans A:AY).
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")

# Columns whose textual values could be misread as numbers/dates by the
# COM layer when written back (e.g. "9", or date-like strings) and so
# must be forced to be stored as plain text rather than being
# reinterpreted as a different type.
$textForceCols = @("I","Y","Z","AA","AB")

function Get-RowValues($rowNum) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$rowNum").Value2
    }
    return $vals
}

function Set-Cell($col, $rowNum, $value) {
    $cell = $ws.Range("$col$rowNum")

    if ($null -eq $value) {
        # Source cell did not exist at all - make sure the target doesn't
        # carry over any stale content either.
        $cell.ClearContents()
        return
    }

    if ([string]::IsNullOrEmpty($value) -and ($value -isnot [bool])) {
        # Source cell existed but was blank - keep a present-but-empty
        # cell behind rather than leaving no cell at all.
        $cell.NumberFormat = "@"
        $cell.Value2 = ""
        $cell.ClearFormats()
        return
    }

    if ($textForceCols -contains $col) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $value
        $cell.ClearFormats()
    } else {
        $cell.Value2 = $value
    }
}

# Snapshot both rows before overwriting anything.
$row2 = Get-RowValues 2
$row3 = Get-RowValues 3

foreach ($col in $cols) {
    Set-Cell $col 2 $row3[$col]
}
foreach ($col in $cols) {
    Set-Cell $col 3 $row2[$col]
}
